# Path To Graduation Y.xlsx - genericize the semester-header placeholders.
#
# The sheet has five "semester header" rows (3, 11, 19, 27, 35), each with
# three labeled mini-tables: a Fall/Spring/Summer header in columns A/C/E.
# They used to be hard-coded to specific academic years (e.g. "Fall 2022",
# "Fall 2023", ... "Fall 2026"), which made it look like courses (like
# CPSC 4000) were scheduled for one specific, fixed term. Replace all of
# them with generic "Fall 20__" / "Spring 20__" / "Summer 20__" placeholders
# so the template isn't tied to (and doesn't misrepresent) a specific year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRows = @(3, 11, 19, 27, 35)

foreach ($row in $headerRows) {
    $ws.Range("A$row").Value = "Fall 20__"
    $ws.Range("C$row").Value = "Spring 20__"
    $ws.Range("E$row").Value = "Summer 20__"
}

# Restore the cursor/selection to where the author left it.
$ws.Range("E35").Select()
